$d = $word.ActiveDocument

# 1. Title replacements (appears twice: main heading + bold text near end)
$d.Content.Find.Execute("Play Da Vinci Extreme Slot for Free | Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Da Vinci Extreme Free | Review", 2)

# 2. Tumbling Reels bullet
$d.Content.Find.Execute("Tumbling Reels format for consecutive wins and free respins", $true, $false, $false, $false, $false, $true, 1, $false, "Tumbling Reels format allows for consecutive wins and respins", 2)

# 3. Free spins bullet
$d.Content.Find.Execute("Free spins bonus feature with up to 15 extra spins available", $true, $false, $false, $false, $false, $true, 1, $false, "Free spins bonus feature with the potential to stack up to 300 free spins", 2)

# 4. Brilliant graphics bullet
$d.Content.Find.Execute("Brilliant graphics that transport players to the world of Leonardo da Vinci", $true, $false, $false, $false, $false, $true, 1, $false, "Brilliant graphics that transport players to the world of Da Vinci", 2)

# 5. Autoplay features bullet
$d.Content.Find.Execute("Autoplay features and advanced settings for customization", $true, $false, $false, $false, $false, $true, 1, $false, "Autoplay features allow for customization of gameplay experience", 2)

# 6/7. Swap the two "What we don't like" bullets: the old first bullet's text
# becomes the (reworded) old second bullet's text, and vice versa.
$d.Content.Find.Execute("Golden wild symbol only appears on reels 2-4", $true, $false, $false, $false, $false, $true, 1, $false, "~~SWAP_PLACEHOLDER~~", 2)
$d.Content.Find.Execute("Basic graphics and music", $true, $false, $false, $false, $false, $true, 1, $false, "Golden wild symbol only appears on reels 2-4", 2)
$d.Content.Find.Execute("~~SWAP_PLACEHOLDER~~", $true, $false, $false, $false, $false, $true, 1, $false, "Graphics and music are basic", 2)

# 8. Bold title again (already covered by step 1's replace-all since both Heading1 and bold text match)

# 9. Italic summary text
$d.Content.Find.Execute("Read our unbiased review of Da Vinci Extreme and play this 5-reel online slot machine for free. Enjoy Tumbling Reels and a free spins bonus feature.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Da Vinci Extreme slot game and play it for free.", 2)
